{"js": "// Replace the division problems in the table with their updated values,\n// matching the author's edit (see commit diff). Each old value is unique\n// in the document, so we can safely search-and-replace each one in turn.\nconst replacements = [\n  [\"426\u00f74=\", \"401\u00f73=\"],\n  [\"561\u00f79=\", \"566\u00f72=\"],\n  [\"386\u00f79=\", \"855\u00f74=\"],\n  [\"483\u00f72=\", \"126\u00f72=\"],\n  [\"369\u00f79=\", \"600\u00f75=\"],\n  [\"755\u00f72=\", \"943\u00f76=\"],\n  [\"213\u00f76=\", \"498\u00f78=\"],\n  [\"892\u00f77=\", \"157\u00f72=\"],\n  [\"186\u00f79=\", \"926\u00f76=\"],\n  [\"904\u00f74=\", \"821\u00f72=\"],\n  [\"514\u00f73=\", \"964\u00f78=\"],\n  [\"608\u00f73=\", \"175\u00f72=\"],\n  [\"794\u00f75=\", \"439\u00f79=\"],\n  [\"920\u00f72=\", \"680\u00f76=\"],\n  [\"709\u00f73=\", \"260\u00f75=\"],\n  [\"504\u00f76=\", \"448\u00f79=\"],\n  [\"565\u00f79=\", \"457\u00f77=\"],\n  [\"876\u00f72=\", \"796\u00f79=\"],\n  [\"962\u00f74=\", \"198\u00f77=\"],\n  [\"451\u00f73=\", \"125\u00f75=\"],\n  [\"990\u00f76=\", \"604\u00f77=\"],\n  [\"428\u00f74=\", \"977\u00f77=\"],\n  [\"972\u00f73=\", \"921\u00f75=\"],\n  [\"944\u00f74=\", \"769\u00f76=\"],\n  [\"752\u00f76=\", \"188\u00f74=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const result of results.items) {\n    result.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the division problems in the table with their updated values,\n# matching the author's edit (see commit diff). Each old value is unique\n# in the document, so Find/Replace for each pair is unambiguous.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"426\u00f74=\", \"401\u00f73=\"),\n    @(\"561\u00f79=\", \"566\u00f72=\"),\n    @(\"386\u00f79=\", \"855\u00f74=\"),\n    @(\"483\u00f72=\", \"126\u00f72=\"),\n    @(\"369\u00f79=\", \"600\u00f75=\"),\n    @(\"755\u00f72=\", \"943\u00f76=\"),\n    @(\"213\u00f76=\", \"498\u00f78=\"),\n    @(\"892\u00f77=\", \"157\u00f72=\"),\n    @(\"186\u00f79=\", \"926\u00f76=\"),\n    @(\"904\u00f74=\", \"821\u00f72=\"),\n    @(\"514\u00f73=\", \"964\u00f78=\"),\n    @(\"608\u00f73=\", \"175\u00f72=\"),\n    @(\"794\u00f75=\", \"439\u00f79=\"),\n    @(\"920\u00f72=\", \"680\u00f76=\"),\n    @(\"709\u00f73=\", \"260\u00f75=\"),\n    @(\"504\u00f76=\", \"448\u00f79=\"),\n    @(\"565\u00f79=\", \"457\u00f77=\"),\n    @(\"876\u00f72=\", \"796\u00f79=\"),\n    @(\"962\u00f74=\", \"198\u00f77=\"),\n    @(\"451\u00f73=\", \"125\u00f75=\"),\n    @(\"990\u00f76=\", \"604\u00f77=\"),\n    @(\"428\u00f74=\", \"977\u00f77=\"),\n    @(\"972\u00f73=\", \"921\u00f75=\"),\n    @(\"944\u00f74=\", \"769\u00f76=\"),\n    @(\"752\u00f76=\", \"188\u00f74=\")\n)\n\nforeach ($pair in $pairs) {\n    $old = $pair[0]\n    $new = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $old\n    $find.Replacement.Text = $new\n    $find.Forward = $true\n    $find.Wrap = 1\n\n    # wdFindContinue=1, wdReplaceAll=2\n    $find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null\n}\n"}
